# Trade #23 closed at 2026-02-16 21:25:47 - momentum DOWN +0.000%
#
# This script:
#  1) Closes trade #4 (momentum strategy) on the "momentum" and mirrors it
#     into "All Trades" (row 5).
#  2) Opens a brand-new trade #23 on the "momentum" sheet (row 6).
#  3) Refreshes the aggregate stats on "Summary" (row 2 + new row 4) and
#     "Comparison" (new row 3).

$wb = $excel.ActiveWorkbook

# Helper: write a value as literal text (prevents Excel from silently
# re-interpreting strings that look like numbers/percents/dates, e.g.
# "50.0%", "+0.6926%", "2026-02-16", "0.00", "inf" ...) while leaving the
# cell's style pointing back at the default "Normal" style (no stray
# numFmt left behind on the cell).
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Helper: write a genuinely empty (but Text-typed) cell, mirroring the
# workbook's own convention of an empty inline string placeholder
# (<c t="inlineStr"/>) for "no value yet" fields such as Exit Price /
# Exit Reason on a still-OPEN trade. A bare "" assignment gets optimised
# away entirely by the engine (cell stays absent / falls back to the
# generic Number type), but a lone quote character is recognised as an
# empty, quote-prefixed string -- resetting the style afterwards drops
# the quote-prefix formatting so the cell is indistinguishable from the
# source file's own blank placeholders.
function Set-EmptyText {
    param($range)
    $range.Value = "'"
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet "momentum": close trade #4 (row 2) and append new trade #23 (row 6)
# ---------------------------------------------------------------------
$wsMom = $wb.Worksheets.Item("momentum")

# Widen a few columns (G, I, M) to match the new content.
$wsMom.Columns.Item(7).ColumnWidth = 14 - 0.875
$wsMom.Columns.Item(9).ColumnWidth = 8 - 0.875
$wsMom.Columns.Item(13).ColumnWidth = 16 - 0.875

# Row 2: trade #4 goes from OPEN -> CLOSED
$wsMom.Range("G2").Value = 68979.847029
$wsMom.Range("H2").Value = "CLOSED"
$wsMom.Range("I2").Value = 0.6926
$wsMom.Range("J2").Value = 6.93
Set-TextValue $wsMom.Range("M2") "time_exit_5min"
$wsMom.Range("N2").Value = 5

# Row 6: brand-new trade #23, still OPEN
$wsMom.Range("A6").Value = 23
Set-TextValue $wsMom.Range("B6") "2026-02-16"
Set-TextValue $wsMom.Range("C6") "21:25:47"
Set-TextValue $wsMom.Range("D6") "momentum"
Set-TextValue $wsMom.Range("E6") "DOWN"
$wsMom.Range("F6").Value = 69090.855
Set-EmptyText $wsMom.Range("G6")
$wsMom.Range("H6").Value = "OPEN"
$wsMom.Range("I6").Value = 0
$wsMom.Range("J6").Value = 0
$wsMom.Range("K6").Value = 0.9
Set-TextValue $wsMom.Range("L6") "Downward momentum: -0.280% over 10 samples"
Set-EmptyText $wsMom.Range("M6")
$wsMom.Range("N6").Value = 0

# ---------------------------------------------------------------------
# Sheet "All Trades": append the now-closed trade #4 as row 5
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Columns.Item(12).ColumnWidth = 44 - 0.875

$wsAll.Range("A5").Value = 4
Set-TextValue $wsAll.Range("B5") "2026-02-16"
Set-TextValue $wsAll.Range("C5") "21:20:43"
Set-TextValue $wsAll.Range("D5") "momentum"
Set-TextValue $wsAll.Range("E5") "DOWN"
$wsAll.Range("F5").Value = 69460.925
$wsAll.Range("G5").Value = 68979.847029
$wsAll.Range("H5").Value = "CLOSED"
$wsAll.Range("I5").Value = 0.6926
$wsAll.Range("J5").Value = 6.93
$wsAll.Range("K5").Value = 0.9
Set-TextValue $wsAll.Range("L5") "Downward momentum: -0.208% over 10 samples"
Set-TextValue $wsAll.Range("M5") "time_exit_5min"
$wsAll.Range("N5").Value = 5

# ---------------------------------------------------------------------
# Sheet "Summary": refresh OVERALL row + add new STRATEGY row for momentum
# ---------------------------------------------------------------------
$wsSum = $wb.Worksheets.Item("Summary")

$wsSum.Range("C2").Value = 4
Set-TextValue $wsSum.Range("D2") "50.0%"
Set-TextValue $wsSum.Range("E2") "+0.6541%"
Set-TextValue $wsSum.Range("F2") "+0.1635%"

Set-TextValue $wsSum.Range("A4") "STRATEGY"
Set-TextValue $wsSum.Range("B4") "momentum"
$wsSum.Range("C4").Value = 4
Set-TextValue $wsSum.Range("D4") "25.0%"
Set-TextValue $wsSum.Range("E4") "+0.6926%"
Set-TextValue $wsSum.Range("F4") "+0.1731%"

# ---------------------------------------------------------------------
# Sheet "Comparison": add new row for momentum strategy
# ---------------------------------------------------------------------
$wsCmp = $wb.Worksheets.Item("Comparison")

Set-TextValue $wsCmp.Range("A3") "momentum"
$wsCmp.Range("B3").Value = 4
Set-TextValue $wsCmp.Range("C3") "25.0%"
Set-TextValue $wsCmp.Range("D3") "inf"
Set-TextValue $wsCmp.Range("E3") "+0.6926%"
Set-TextValue $wsCmp.Range("F3") "+0.0000%"
Set-TextValue $wsCmp.Range("G3") "0.00"
Set-TextValue $wsCmp.Range("H3") "+0.0000%"

Write-Host "edit applied"
